$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2025-10-05 18:30:13", 5078),
    @("2025-10-05 18:30:16", 704),
    @("2025-10-05 18:30:20", 8477),
    @("2025-10-05 18:30:24", 1654),
    @("2025-10-05 18:30:27", 7698),
    @("2025-10-05 18:30:31", 2339),
    @("2025-10-05 18:30:35", 9581),
    @("2025-10-05 18:30:38", 9874),
    @("2025-10-05 18:30:42", 4600),
    @("2025-10-05 18:30:45", 2989),
    @("2025-10-05 18:30:49", 2648),
    @("2025-10-05 18:30:53", 4390),
    @("2025-10-05 18:30:56", 8799),
    @("2025-10-05 18:31:00", 3997),
    @("2025-10-05 18:31:04", 912),
    @("2025-10-05 18:36:38", 8503),
    @("2025-10-05 18:36:42", 7744),
    @("2025-10-05 18:36:45", 7388),
    @("2025-10-05 18:36:49", 7690),
    @("2025-10-05 18:36:53", 9189),
    @("2025-10-05 18:36:56", 1373),
    @("2025-10-05 18:37:00", 4024),
    @("2025-10-05 18:37:04", 9498),
    @("2025-10-05 18:37:07", 6539),
    @("2025-10-05 18:37:11", 3381),
    @("2025-10-05 18:37:15", 2398),
    @("2025-10-05 18:37:18", 5299),
    @("2025-10-05 18:37:22", 6327),
    @("2025-10-05 18:37:25", 3293),
    @("2025-10-05 18:37:29", 1890),
    @("2025-10-05 18:37:33", 6090),
    @("2025-10-05 18:37:36", 3860),
    @("2025-10-05 18:37:40", 6387),
    @("2025-10-05 18:37:44", 2190),
    @("2025-10-05 18:37:47", 2625),
    @("2025-10-05 18:37:51", 489),
    @("2025-10-05 18:37:54", 7416),
    @("2025-10-05 18:37:58", 4484),
    @("2025-10-05 18:38:02", 3765),
    @("2025-10-05 18:38:05", 9606),
    @("2025-10-05 18:38:09", 9401),
    @("2025-10-05 18:38:13", 1321),
    @("2025-10-05 18:38:16", 4658),
    @("2025-10-05 18:38:20", 8473),
    @("2025-10-05 18:38:24", 252),
    @("2025-10-05 18:38:27", 775),
    @("2025-10-05 18:38:31", 7497),
    @("2025-10-05 18:38:34", 4080),
    @("2025-10-05 18:38:38", 8536),
    @("2025-10-05 18:38:42", 4126),
    @("2025-10-05 18:38:45", 9741),
    @("2025-10-05 18:38:49", 2994),
    @("2025-10-05 18:38:53", 4568),
    @("2025-10-05 18:38:56", 5569),
    @("2025-10-05 18:39:00", 239),
    @("2025-10-05 18:39:04", 3785),
    @("2025-10-05 18:39:07", 8587),
    @("2025-10-05 18:39:11", 7762),
    @("2025-10-05 18:39:14", 5898),
    @("2025-10-05 18:39:18", 3503),
    @("2025-10-05 18:39:22", 9516),
    @("2025-10-05 18:39:25", 5988),
    @("2025-10-05 18:39:29", 386),
    @("2025-10-05 18:39:33", 7810),
    @("2025-10-05 18:39:36", 6606),
    @("2025-10-05 18:39:40", 3592),
    @("2025-10-05 18:39:44", 5277),
    @("2025-10-05 18:39:47", 1269),
    @("2025-10-05 18:39:51", 5855),
    @("2025-10-05 18:39:55", 7127),
    @("2025-10-05 18:39:58", 7679),
    @("2025-10-05 18:40:02", 2328),
    @("2025-10-05 18:40:05", 8589),
    @("2025-10-05 18:40:09", 6267),
    @("2025-10-05 18:40:13", 9458),
    @("2025-10-05 18:40:16", 8112),
    @("2025-10-05 18:40:20", 9246),
    @("2025-10-05 18:40:24", 9552),
    @("2025-10-05 18:40:27", 9411),
    @("2025-10-05 18:40:31", 5089),
    @("2025-10-05 18:40:35", 9866),
    @("2025-10-05 18:40:38", 3231),
    @("2025-10-05 18:40:42", 9778),
    @("2025-10-05 18:40:46", 5113),
    @("2025-10-05 18:40:49", 4541),
    @("2025-10-05 18:40:53", 2064),
    @("2025-10-05 18:40:57", 3487),
    @("2025-10-05 18:41:00", 3224),
    @("2025-10-05 18:41:04", 5423),
    @("2025-10-05 18:41:08", 8311),
    @("2025-10-05 18:41:11", 9182),
    @("2025-10-05 18:41:15", 1673),
    @("2025-10-05 18:41:18", 1741),
    @("2025-10-05 18:41:22", 3736),
    @("2025-10-05 18:41:26", 2034),
    @("2025-10-05 18:41:29", 1077),
    @("2025-10-05 18:41:33", 6076),
    @("2025-10-05 18:41:37", 6786),
    @("2025-10-05 18:41:40", 2604),
    @("2025-10-05 18:41:44", 4164),
    @("2025-10-05 18:41:48", 7202),
    @("2025-10-05 18:41:51", 8341),
    @("2025-10-05 18:41:55", 8893),
    @("2025-10-05 18:41:59", 832),
    @("2025-10-05 18:42:02", 963),
    @("2025-10-05 18:42:06", 2683),
    @("2025-10-05 18:42:10", 6540),
    @("2025-10-05 18:42:13", 2307),
    @("2025-10-05 18:42:17", 5035),
    @("2025-10-05 18:42:21", 4485),
    @("2025-10-05 18:42:24", 3266),
    @("2025-10-05 18:42:28", 7454),
    @("2025-10-05 18:42:32", 39),
    @("2025-10-05 18:42:35", 3466),
    @("2025-10-05 18:42:39", 5612)
)

$startRow = 102
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

Write-Output "Appended $($data.Count) rows"